$d = $word.ActiveDocument

# Replace "/5kmp/2BS" with "/5kp/2BS" (kmp -> kp) throughout the document body
$d.Content.Find.Execute("/5kmp/2BS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "/5kp/2BS", 2)
